$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.282.14'
$ws.Range('E2').Value = '  +0.59%  '
$ws.Range('D3').Value = '1.664.03'
$ws.Range('E3').Value = '  +0.67%  '
$ws.Range('E4').Value = '  +0.74%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.63'
$ws.Range('E5').Value = '  +0.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5326'
$ws.Range('E6').Value = '  +0.71%  '
$ws.Range('E7').Value = '  +0.69%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2642'
$ws.Range('E8').Value = '  +1.28%  '
$ws.Range('E9').Value = '  +0.52%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.56'
$ws.Range('E10').Value = '  +0.78%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07848'
$ws.Range('E11').Value = '  +1.25%  '
$ws.Range('E12').Value = '  +1.46%  '
$ws.Range('D13').Value = '1.665.79'
$ws.Range('E13').Value = '  -0.45%  '
$ws.Range('D14').Value = '1.893.27'
$ws.Range('E14').Value = '  +0.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5532'
$ws.Range('E15').Value = '  +1.26%  '
$ws.Range('D16').Value = '0.0₅8197'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.71'
$ws.Range('E17').Value = '  +0.65%  '
$ws.Range('E18').Value = '  +0.73%  '
$ws.Range('E19').Value = '  +2.86%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '192.64'
$ws.Range('E20').Value = '  -0.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.21'
$ws.Range('E21').Value = '  +1.82%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.059'
$ws.Range('E22').Value = '  +1.14%  '
$ws.Range('E23').Value = '  +0.66%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '145.12'
$ws.Range('E24').Value = '  +3.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1226'
$ws.Range('E25').Value = '  -1.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.253'
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.14'
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('E28').Value = '  +2.76%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.05851'
$ws.Range('E29').Value = '  -1.41%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.278'
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.586'
$ws.Range('E31').Value = '  +2.13%  '
$ws.Range('E32').Value = '  +2.22%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.620'
$ws.Range('E33').Value = '  +4.50%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9595'
$ws.Range('E34').Value = '  +1.50%  '
$ws.Range('B35').Value = 'MXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.821'
$ws.Range('E35').Value = '  +2.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.429'
$ws.Range('E36').Value = '  +0.68%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5820'
$ws.Range('E37').Value = '  +3.33%  '
$ws.Range('E38').Value = '  +0.53%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.890'
$ws.Range('E39').Value = '  +0.54%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8545'
$ws.Range('E40').Value = '  +0.96%  '
$ws.Range('E41').Value = '  +0.70%  '
$ws.Range('D42').Value = '1.047.81'
$ws.Range('E42').Value = '  +3.87%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '104.54'
$ws.Range('E43').Value = '  +3.68%  '
$ws.Range('D44').Value = '1.806.01'
$ws.Range('E44').Value = '  +0.58%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '57.39'
$ws.Range('E45').Value = '  +1.09%  '
$ws.Range('D46').Value = '0.0₈108'
$ws.Range('E46').Value = '  +2.50%  '
$ws.Range('E47').Value = '  +0.44%  '
$ws.Range('E48').Value = '  +1.98%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.945'
$ws.Range('E49').Value = '  +2.11%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05164'
$ws.Range('E50').Value = '  +0.25%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.448'
$ws.Range('E51').Value = '  -1.63%  '
